$d = $word.ActiveDocument

# Update the date heading at the top of the document
$d.Content.Find.Execute("2025-07-04 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-05 Saturday", 2) | Out-Null

# Update the multiplication answers in the table, cell by cell so that
# duplicate values (e.g. the two "683x5=3415" cells) are each replaced
# with their own distinct replacement rather than a document-wide swap.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("517×4=2068", $true, $false, $false, $false, $false, $true, 0, $false, "782×8=6256", 1) | Out-Null

$cell = $t.Cell(1, 2)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("986×5=4930", $true, $false, $false, $false, $false, $true, 0, $false, "542×4=2168", 1) | Out-Null

$cell = $t.Cell(1, 3)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("357×5=1785", $true, $false, $false, $false, $false, $true, 0, $false, "775×9=6975", 1) | Out-Null

$cell = $t.Cell(1, 4)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("465×7=3255", $true, $false, $false, $false, $false, $true, 0, $false, "988×8=7904", 1) | Out-Null

$cell = $t.Cell(1, 5)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("372×4=1488", $true, $false, $false, $false, $false, $true, 0, $false, "543×2=1086", 1) | Out-Null

$cell = $t.Cell(5, 1)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("687×4=2748", $true, $false, $false, $false, $false, $true, 0, $false, "292×9=2628", 1) | Out-Null

$cell = $t.Cell(5, 2)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("716×7=5012", $true, $false, $false, $false, $false, $true, 0, $false, "489×6=2934", 1) | Out-Null

$cell = $t.Cell(5, 3)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("551×7=3857", $true, $false, $false, $false, $false, $true, 0, $false, "475×9=4275", 1) | Out-Null

$cell = $t.Cell(5, 4)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("272×3=816", $true, $false, $false, $false, $false, $true, 0, $false, "935×9=8415", 1) | Out-Null

$cell = $t.Cell(5, 5)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("903×7=6321", $true, $false, $false, $false, $false, $true, 0, $false, "605×7=4235", 1) | Out-Null

$cell = $t.Cell(10, 1)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("557×2=1114", $true, $false, $false, $false, $false, $true, 0, $false, "861×9=7749", 1) | Out-Null

$cell = $t.Cell(10, 2)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("332×7=2324", $true, $false, $false, $false, $false, $true, 0, $false, "538×4=2152", 1) | Out-Null

$cell = $t.Cell(10, 3)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("891×2=1782", $true, $false, $false, $false, $false, $true, 0, $false, "627×2=1254", 1) | Out-Null

$cell = $t.Cell(10, 4)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("715×6=4290", $true, $false, $false, $false, $false, $true, 0, $false, "472×5=2360", 1) | Out-Null

$cell = $t.Cell(10, 5)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("514×7=3598", $true, $false, $false, $false, $false, $true, 0, $false, "919×4=3676", 1) | Out-Null

$cell = $t.Cell(15, 1)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("223×7=1561", $true, $false, $false, $false, $false, $true, 0, $false, "849×7=5943", 1) | Out-Null

$cell = $t.Cell(15, 2)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("494×6=2964", $true, $false, $false, $false, $false, $true, 0, $false, "278×7=1946", 1) | Out-Null

$cell = $t.Cell(15, 3)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("683×5=3415", $true, $false, $false, $false, $false, $true, 0, $false, "509×3=1527", 1) | Out-Null

$cell = $t.Cell(15, 4)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("393×7=2751", $true, $false, $false, $false, $false, $true, 0, $false, "610×3=1830", 1) | Out-Null

$cell = $t.Cell(15, 5)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("685×6=4110", $true, $false, $false, $false, $false, $true, 0, $false, "178×7=1246", 1) | Out-Null

$cell = $t.Cell(20, 1)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("683×5=3415", $true, $false, $false, $false, $false, $true, 0, $false, "536×9=4824", 1) | Out-Null

$cell = $t.Cell(20, 2)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("986×3=2958", $true, $false, $false, $false, $false, $true, 0, $false, "560×6=3360", 1) | Out-Null

$cell = $t.Cell(20, 3)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("319×4=1276", $true, $false, $false, $false, $false, $true, 0, $false, "636×4=2544", 1) | Out-Null

$cell = $t.Cell(20, 4)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("464×5=2320", $true, $false, $false, $false, $false, $true, 0, $false, "266×2=532", 1) | Out-Null

$cell = $t.Cell(20, 5)
$cellStart = $cell.Range.Start
$cellEnd = $cell.Range.End
$rng = $d.Range($cellStart, $cellEnd)
$rng.Find.Execute("967×3=2901", $true, $false, $false, $false, $false, $true, 0, $false, "656×9=5904", 1) | Out-Null
